# Add a new "canonical SMILES" column (D) to the microstate list sheet.
# Column D mirrors column C ("canonical isomeric SMILES") for every data
# row, since these microstates have no additional stereo information beyond
# what's already captured canonically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in D2
$ws.Range("D2").Value = "canonical SMILES"

# Fill D3:D22 with the same canonical SMILES values already present in C3:C22
for ($row = 3; $row -le 22; $row++) {
    $source = $ws.Cells.Item($row, 3)
    $target = $ws.Cells.Item($row, 4)
    $target.Value = $source.Value()
}

# Match the column width used for the new column
$ws.Columns.Item(4).ColumnWidth = 36.85546875

Write-Host "Added canonical SMILES column (D) with $(22 - 3 + 1) data rows."
